# Update the single-column results table in place.
# Each row of the table holds one value in its sole cell; we address
# cells directly by (row, column) to avoid ambiguity from duplicate
# values (e.g. several rows contain "0.08227" or "100.0").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements.
$t.Cell(1, 1).Range.Text  = "0M"        # was 96.74
$t.Cell(2, 1).Range.Text  = "0M"        # was 7.34
$t.Cell(3, 1).Range.Text  = "0M"        # was 224
$t.Cell(4, 1).Range.Text  = "1481"      # was 1278
$t.Cell(6, 1).Range.Text  = "0.08457"   # was 0.08227
$t.Cell(7, 1).Range.Text  = "0.02291"   # was 0.02248
$t.Cell(8, 1).Range.Text  = "0.00370"   # was 0.00271
$t.Cell(12, 1).Range.Text = "7.33986"   # was 1.19122

# These cells previously held a whole tab-separated stats row; they
# collapse down to just the lead value (which moved from the rows
# replaced with "0M" above).
$t.Cell(44, 1).Range.Text = "96.74"
$t.Cell(45, 1).Range.Text = "7.34"
$t.Cell(46, 1).Range.Text = "224"
